$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("test_suite")
$ws2 = $wb.Worksheets.Item("AddCustomerTest")
$ws3 = $wb.Worksheets.Item("OpenAccountTest")

# --- test_suite: flip AddCustomerTest and OpenAccountTest Runmode from N to Y ---
$ws1.Range("B3").Value = "Y"
$ws1.Range("B4").Value = "Y"
$ws1.Columns("A").ColumnWidth = 22.2

# --- AddCustomerTest: insert a "runmode" column marking which row to skip ---
$ws2.Columns("A").Insert()
$ws2.Range("A1").Value = "runmode"
$ws2.Range("A2").Value = "N"
$ws2.Range("A3").Value = "Y"

# --- OpenAccountTest: use Rahul Jadhwani (the still-enabled customer) instead of Deepender Singh ---
$ws3.Range("A2").Value = "Rahul Jadhwani"

# --- selection / active-sheet bookkeeping (AddCustomerTest ends up the active tab) ---
$ws1.Range("B9").Select() | Out-Null
$excel.ActiveWindow.Zoom = 85
$ws3.Range("A3").Select() | Out-Null
$ws2.Range("A3").Select() | Out-Null
